$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 7

# Make Sheet1 the active sheet and reset the view / selection
$ws.Activate()
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
